# Actualización automática 2025-07-03 12:35:11
# Update figures on the "CUMPLIMIENTO MENSUAL" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column F width: 26 -> 25
# (Excel's ColumnWidth setter re-quantizes to the default font's pixel grid,
#  so request a value that rounds-trips to a stored width of exactly 25.)
$ws.Columns.Item(6).ColumnWidth = 24.17

# Row 3 - 240X80 PORCELANATO
$ws.Range("D3").Value = 442.27
$ws.Range("E3").Value = 3725.80156573679
$ws.Range("F3").Value = 0.1061090226078736

# Row 4 - FREGADEROS DE COCINA
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 513.831046659336
$ws.Range("F4").Value = 0

# Row 6 - GRIFERIAS
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 106.82
$ws.Range("F6").Value = 0

# Row 7 - INODOROS
$ws.Range("C7").Value = 2400
$ws.Range("D7").Value = 603
$ws.Range("E7").Value = 1797
$ws.Range("F7").Value = 0.25125

# Row 8 - LAVABOS
$ws.Range("D8").Value = 28.8
$ws.Range("E8").Value = 596.2
$ws.Range("F8").Value = 0.04608

# Row 13 - PANELES PU
$ws.Range("C13").Value = 130
$ws.Range("E13").Value = 130

# Row 14 - PANELES PVC
$ws.Range("C14").Value = 240
$ws.Range("D14").Value = 722.97
$ws.Range("E14").Value = -482.97
$ws.Range("F14").Value = 3.012375

# Row 15 - PIEDRA SINTERIZADA
$ws.Range("D15").Value = 86.56999999999999
$ws.Range("E15").Value = 7378.43
$ws.Range("F15").Value = 0.01159678499665104

# Row 16 - PORCELANATO
$ws.Range("C16").Value = 44266.24
$ws.Range("D16").Value = 7291
$ws.Range("E16").Value = 36975.24
$ws.Range("F16").Value = 0.1647079128473528

# Row 17 - PUERTAS DE SEGURIDAD
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 342
$ws.Range("F17").Value = 0

# Row 18 - SAL SOLUBLE
$ws.Range("D18").Value = 1070.53
$ws.Range("E18").Value = 1729.47
$ws.Range("F18").Value = 0.3823321428571428

# Row 19 - TOTAL
$ws.Range("C19").Value = 65377.99762291768
$ws.Range("D19").Value = 10245.14
$ws.Range("E19").Value = 55132.85762291768
$ws.Range("F19").Value = 0.1567062371516967
